# Rotate the data rows 3-6 on the active worksheet.
# New row 3 = old row 5
# New row 4 = old row 3
# New row 5 = old row 6
# New row 6 = old row 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 6
$lastCol = 51   # column AY

# Columns that hold date/time text (Startdatum, Starttid, Slutdatum, Sluttid).
# These must stay plain text instead of being auto-parsed into date/time
# serial numbers when re-assigned.
$textCols = @(25, 26, 27, 28)

# Capture the original values for rows 3-6, columns A:AY
$rowData = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $values = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $values += $ws.Cells.Item($r, $c).Value2
    }
    $rowData[$r] = $values
}

# Mapping: destination row -> source row
$mapping = @{
    3 = 5
    4 = 3
    5 = 6
    6 = 4
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcValues = $rowData[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($destRow, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value2 = $srcValues[$c - 1]
    }
}
